$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 309473
$ws.Range("D2").Value = 394789073
$ws.Range("C3").Value = 252
$ws.Range("D3").Value = 300350
$ws.Range("C4").Value = 309
$ws.Range("D4").Value = 441918
$ws.Range("C8").Value = 831
$ws.Range("D8").Value = 1224107
$ws.Range("C10").Value = 114254
$ws.Range("D10").Value = 167482684
$ws.Range("C11").Value = 140
$ws.Range("D11").Value = 207013
$ws.Range("C12").Value = 57400
$ws.Range("D12").Value = 82868470
$ws.Range("C16").Value = 3818
$ws.Range("D16").Value = 5421792
$ws.Range("C19").Value = 68
$ws.Range("D19").Value = 98535
$ws.Range("C20").Value = 6162
$ws.Range("D20").Value = 8610222
$ws.Range("C22").Value = 75148
$ws.Range("D22").Value = 93902311
$ws.Range("C23").Value = 51
$ws.Range("D23").Value = 66709
$ws.Range("C27").Value = 284
$ws.Range("D27").Value = 407147
$ws.Range("C28").Value = 31771
$ws.Range("D28").Value = 46522351
$ws.Range("C30").Value = 11143
$ws.Range("D30").Value = 16045328
$ws.Range("C33").Value = 1490
$ws.Range("D33").Value = 2094294
$ws.Range("C35").Value = 1692
$ws.Range("D35").Value = 2383549
$ws.Range("C36").Value = 94611
$ws.Range("D36").Value = 119292982
$ws.Range("C37").Value = 64
$ws.Range("D37").Value = 74727
$ws.Range("C42").Value = 888
$ws.Range("D42").Value = 1308272
$ws.Range("C44").Value = 43573
$ws.Range("D44").Value = 63883540
$ws.Range("C46").Value = 8889
$ws.Range("D46").Value = 12763932
$ws.Range("C48").Value = 1359
$ws.Range("D48").Value = 1885019
$ws.Range("C51").Value = 2146
$ws.Range("D51").Value = 2990888
$ws.Range("C52").Value = 66919
$ws.Range("D52").Value = 84020499
$ws.Range("C58").Value = 27521
$ws.Range("D58").Value = 40366157
$ws.Range("C61").Value = 10731
$ws.Range("D61").Value = 15521011
$ws.Range("C63").Value = 1323
$ws.Range("D63").Value = 1846789
$ws.Range("C67").Value = 1379
$ws.Range("D67").Value = 1927403
$ws.Range("C69").Value = 19780
$ws.Range("D69").Value = 25906320
$ws.Range("C70").Value = 31
$ws.Range("D70").Value = 41126
$ws.Range("C72").Value = 55
$ws.Range("D72").Value = 80073
$ws.Range("C73").Value = 7311
$ws.Range("D73").Value = 10703227
$ws.Range("C75").Value = 4925
$ws.Range("D75").Value = 7153369
$ws.Range("C76").Value = 465
$ws.Range("D76").Value = 656239
$ws.Range("C77").Value = 259
$ws.Range("D77").Value = 363213
$ws.Range("C78").Value = 136395
$ws.Range("D78").Value = 170213634
$ws.Range("C83").Value = 13
$ws.Range("D83").Value = 17438
$ws.Range("C84").Value = 62144
$ws.Range("D84").Value = 91106030
$ws.Range("C87").Value = 28778
$ws.Range("D87").Value = 41657015
$ws.Range("C89").Value = 2611
$ws.Range("D89").Value = 3758692
$ws.Range("C90").Value = 2609
$ws.Range("D90").Value = 3683347
$ws.Range("C91").Value = 30328
$ws.Range("D91").Value = 41098166
$ws.Range("C95").Value = 7521
$ws.Range("D95").Value = 11071620
$ws.Range("C97").Value = 6789
$ws.Range("D97").Value = 9841716
$ws.Range("C99").Value = 483
$ws.Range("D99").Value = 687405
$ws.Range("C100").Value = 452
$ws.Range("D100").Value = 653273
$ws.Range("C101").Value = 7940
$ws.Range("D101").Value = 11019065
$ws.Range("C103").Value = 2021
$ws.Range("D103").Value = 2975415
$ws.Range("C105").Value = 2737
$ws.Range("D105").Value = 3996373
$ws.Range("C108").Value = 138
$ws.Range("D108").Value = 196689
$ws.Range("C109").Value = 137023
$ws.Range("D109").Value = 169528862
$ws.Range("C113").Value = 937
$ws.Range("D113").Value = 1375882
$ws.Range("C115").Value = 51605
$ws.Range("D115").Value = 75672727
$ws.Range("C117").Value = 26138
$ws.Range("D117").Value = 37871454
$ws.Range("C118").Value = 1259
$ws.Range("D118").Value = 1723356
$ws.Range("C121").Value = 2119
$ws.Range("D121").Value = 2973107
$ws.Range("C123").Value = 476606
$ws.Range("D123").Value = 628508024
$ws.Range("C128").Value = 1324
$ws.Range("D128").Value = 1962311
$ws.Range("C129").Value = 30
$ws.Range("D129").Value = 39010
$ws.Range("C130").Value = 200279
$ws.Range("D130").Value = 294512507
$ws.Range("C131").Value = 372
$ws.Range("D131").Value = 554790
$ws.Range("C133").Value = 172270
$ws.Range("D133").Value = 250457545
$ws.Range("C136").Value = 2686
$ws.Range("D136").Value = 3771597
$ws.Range("C138").Value = 5826
$ws.Range("D138").Value = 8232797
$ws.Range("C141").Value = 42695
$ws.Range("D141").Value = 57074279
$ws.Range("C144").Value = 6
$ws.Range("D144").Value = 9000
$ws.Range("C147").Value = 13637
$ws.Range("D147").Value = 20010076
$ws.Range("C148").Value = 3623
$ws.Range("D148").Value = 5227267
$ws.Range("C151").Value = 373
$ws.Range("D151").Value = 536202
$ws.Range("C153").Value = 352
$ws.Range("D153").Value = 495010
$ws.Range("C154").Value = 16651
$ws.Range("D154").Value = 22024026
$ws.Range("C157").Value = 39
$ws.Range("D157").Value = 56906
$ws.Range("C158").Value = 6866
$ws.Range("D158").Value = 9993822
$ws.Range("C160").Value = 4727
$ws.Range("D160").Value = 6808755
$ws.Range("C162").Value = 265
$ws.Range("D162").Value = 365559
$ws.Range("C163").Value = 247
$ws.Range("D163").Value = 353933
$ws.Range("C165").Value = 14071
$ws.Range("D165").Value = 20410134
$ws.Range("C166").Value = 1671
$ws.Range("D166").Value = 2485130
$ws.Range("C167").Value = 224
$ws.Range("D167").Value = 330802
$ws.Range("C171").Value = 85305
$ws.Range("D171").Value = 106773857
$ws.Range("C172").Value = 29
$ws.Range("D172").Value = 36229
$ws.Range("C176").Value = 632
$ws.Range("D176").Value = 931348
$ws.Range("C178").Value = 33156
$ws.Range("D178").Value = 48635526
$ws.Range("C180").Value = 12605
$ws.Range("D180").Value = 18215592
$ws.Range("C182").Value = 1199
$ws.Range("D182").Value = 1678257
$ws.Range("C184").Value = 1540
$ws.Range("D184").Value = 2165709
$ws.Range("C186").Value = 230824
$ws.Range("D186").Value = 287152057
$ws.Range("C188").Value = 163
$ws.Range("D188").Value = 234736
$ws.Range("C192").Value = 859
$ws.Range("D192").Value = 1263497
$ws.Range("C194").Value = 84783
$ws.Range("D194").Value = 124306984
$ws.Range("C195").Value = 92
$ws.Range("D195").Value = 133627
$ws.Range("C197").Value = 32046
$ws.Range("D197").Value = 46127048
$ws.Range("C200").Value = 4888
$ws.Range("D200").Value = 6964365
$ws.Range("C203").Value = 4522
$ws.Range("D203").Value = 6255652
$ws.Range("C206").Value = 254361
$ws.Range("D206").Value = 314916624
$ws.Range("C207").Value = 152
$ws.Range("D207").Value = 166518
$ws.Range("C208").Value = 244
$ws.Range("D208").Value = 349064
$ws.Range("C213").Value = 602
$ws.Range("D213").Value = 876906
$ws.Range("C215").Value = 92676
$ws.Range("D215").Value = 135622863
$ws.Range("C218").Value = 49743
$ws.Range("D218").Value = 71927879
$ws.Range("C221").Value = 4489
$ws.Range("D221").Value = 6298404
$ws.Range("C224").Value = 5313
$ws.Range("D224").Value = 7345090
$ws.Range("C227").Value = 103072
$ws.Range("D227").Value = 129145086
$ws.Range("C229").Value = 71
$ws.Range("D229").Value = 101945
$ws.Range("C232").Value = 558
$ws.Range("D232").Value = 816339
$ws.Range("C234").Value = 48498
$ws.Range("D234").Value = 71066993
$ws.Range("C236").Value = 11987
$ws.Range("D236").Value = 17235823
$ws.Range("C238").Value = 1837
$ws.Range("D238").Value = 2634659
$ws.Range("C240").Value = 2366
$ws.Range("D240").Value = 3304968
$ws.Range("C241").Value = 248729
$ws.Range("D241").Value = 314298908
$ws.Range("C242").Value = 164
$ws.Range("D242").Value = 204290
$ws.Range("C243").Value = 243
$ws.Range("D243").Value = 348957
$ws.Range("C247").Value = 796
$ws.Range("D247").Value = 1169896
$ws.Range("C249").Value = 93512
$ws.Range("D249").Value = 137060782
$ws.Range("C252").Value = 62766
$ws.Range("D252").Value = 90973646
$ws.Range("C254").Value = 2328
$ws.Range("D254").Value = 3286049
$ws.Range("C257").Value = 4277
$ws.Range("D257").Value = 6005261
